# Insert a new data row at row 96 (pushes existing rows 96-152 down to 97-153)
# and populate it with a new Pomelo "Start Ruby / Especial" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(96).Insert()

$ws.Range("A96").Value = 10
$ws.Range("B96").Value = "Vega Modelo de Temuco"
$ws.Range("C96").Value = "La Araucanía"
$ws.Range("D96").Value = 44488
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100102
$ws.Range("H96").Value = "Cítricos"
$ws.Range("I96").Value = 100102006
$ws.Range("J96").Value = "Pomelo"
$ws.Range("K96").Value = "Start Ruby"
$ws.Range("L96").Value = "Especial"
$ws.Range("M96").Value = 50
$ws.Range("N96").Value = 16000
$ws.Range("O96").Value = 16000
$ws.Range("P96").Value = 16000
$ws.Range("Q96").Value = "$/caja 14 kilos empedrada"
$ws.Range("R96").Value = "Región de O'Higgins"
$ws.Range("S96").Value = 1143
$ws.Range("T96").Value = 14
